$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.403.04'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '1.802.35'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '228.07'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = '0.582'
$ws.Range('E6').Value = '  +4.36%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '35.00'
$ws.Range('E8').Value = '  +6.34%  '
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = '0.0692'
$ws.Range('E10').Value = '  -0.25%  '
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '2.062.66'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.827.02'
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '11.18'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').Value = '0.642'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = '34.390.30'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('D18').Value = '68.96'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '245.38'
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0796'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').Value = '11.50'
$ws.Range('E21').Value = '  +2.13%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').Value = '170.29'
$ws.Range('E24').Value = '  +3.13%  '
$ws.Range('E25').Value = '  +1.99%  '
$ws.Range('D26').Value = '7.62'
$ws.Range('E26').Value = '  +4.88%  '
$ws.Range('E27').Value = '  +2.50%  '
$ws.Range('D28').Value = '16.71'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('D30').Value = '3.99'
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').Value = '0.0528'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').Value = '1.24'
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('D35').Value = '1.395.62'
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('D36').Value = '0.678'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').Value = '83.00'
$ws.Range('E40').Value = '  -2.66%  '
$ws.Range('D41').Value = '2.84'
$ws.Range('E41').Value = '  +2.97%  '
$ws.Range('D42').Value = '0.946'
$ws.Range('E42').Value = '  +1.48%  '
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').Value = '13.55'
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('D46').Value = '0.0509'
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('E47').Value = '  -2.11%  '
$ws.Range('D48').Value = '1.962.80'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '104.40'
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('E51').Value = '  +0.82%  '
